{"js": "// Update the date line and every \"AxB=C\" equation cell in the table.\n// All old values are unique in the document, so a plain body.search()\n// + insertText(\"Replace\") per pair is safe and order-independent.\nconst replacements = [\n  [\"2025-07-11 Friday\", \"2025-07-12 Saturday\"],\n  [\"41\u00d731=1271\", \"14\u00d715=210\"],\n  [\"73\u00d768=4964\", \"55\u00d772=3960\"],\n  [\"19\u00d749=931\", \"55\u00d760=3300\"],\n  [\"36\u00d742=1512\", \"26\u00d776=1976\"],\n  [\"32\u00d723=736\", \"51\u00d737=1887\"],\n  [\"15\u00d720=300\", \"27\u00d780=2160\"],\n  [\"17\u00d732=544\", \"74\u00d795=7030\"],\n  [\"63\u00d714=882\", \"40\u00d712=480\"],\n  [\"17\u00d752=884\", \"95\u00d734=3230\"],\n  [\"15\u00d724=360\", \"84\u00d791=7644\"],\n  [\"48\u00d736=1728\", \"15\u00d756=840\"],\n  [\"24\u00d737=888\", \"62\u00d726=1612\"],\n  [\"38\u00d784=3192\", \"24\u00d733=792\"],\n  [\"38\u00d740=1520\", \"67\u00d724=1608\"],\n  [\"50\u00d742=2100\", \"43\u00d716=688\"],\n  [\"65\u00d754=3510\", \"22\u00d783=1826\"],\n  [\"87\u00d719=1653\", \"90\u00d777=6930\"],\n  [\"70\u00d773=5110\", \"77\u00d743=3311\"],\n  [\"86\u00d715=1290\", \"60\u00d755=3300\"],\n  [\"96\u00d799=9504\", \"65\u00d774=4810\"],\n  [\"92\u00d755=5060\", \"21\u00d772=1512\"],\n  [\"15\u00d723=345\", \"86\u00d790=7740\"],\n  [\"40\u00d734=1360\", \"17\u00d750=850\"],\n  [\"75\u00d752=3900\", \"58\u00d734=1972\"],\n  [\"14\u00d725=350\", \"49\u00d719=931\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"AxB=C\" equation cell in the table.\n# All old values are unique in the document, so a Find/Replace (one\n# match each) per pair is safe and order-independent.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-07-11 Friday\", \"2025-07-12 Saturday\"),\n    @(\"41\u00d731=1271\", \"14\u00d715=210\"),\n    @(\"73\u00d768=4964\", \"55\u00d772=3960\"),\n    @(\"19\u00d749=931\", \"55\u00d760=3300\"),\n    @(\"36\u00d742=1512\", \"26\u00d776=1976\"),\n    @(\"32\u00d723=736\", \"51\u00d737=1887\"),\n    @(\"15\u00d720=300\", \"27\u00d780=2160\"),\n    @(\"17\u00d732=544\", \"74\u00d795=7030\"),\n    @(\"63\u00d714=882\", \"40\u00d712=480\"),\n    @(\"17\u00d752=884\", \"95\u00d734=3230\"),\n    @(\"15\u00d724=360\", \"84\u00d791=7644\"),\n    @(\"48\u00d736=1728\", \"15\u00d756=840\"),\n    @(\"24\u00d737=888\", \"62\u00d726=1612\"),\n    @(\"38\u00d784=3192\", \"24\u00d733=792\"),\n    @(\"38\u00d740=1520\", \"67\u00d724=1608\"),\n    @(\"50\u00d742=2100\", \"43\u00d716=688\"),\n    @(\"65\u00d754=3510\", \"22\u00d783=1826\"),\n    @(\"87\u00d719=1653\", \"90\u00d777=6930\"),\n    @(\"70\u00d773=5110\", \"77\u00d743=3311\"),\n    @(\"86\u00d715=1290\", \"60\u00d755=3300\"),\n    @(\"96\u00d799=9504\", \"65\u00d774=4810\"),\n    @(\"92\u00d755=5060\", \"21\u00d772=1512\"),\n    @(\"15\u00d723=345\", \"86\u00d790=7740\"),\n    @(\"40\u00d734=1360\", \"17\u00d750=850\"),\n    @(\"75\u00d752=3900\", \"58\u00d734=1972\"),\n    @(\"14\u00d725=350\", \"49\u00d719=931\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2) | Out-Null\n}\n\n$d.Save()\n"}
